# Applies the edits described by the commit diff:
#  - merges several runs of unchanged text into single runs (same visible
#    text, same formatting per merged group) across a handful of paragraphs
#  - removes the spell-check proofErr wrapper around "clearcut"
#  - removes a couple of now-redundant empty paragraphs
#  - removes the _GoBack bookmark
#  - adds justify (both) alignment to a couple of paragraphs
#  - repositions two floating pictures

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# 1) "client-server model. " (bold run) / "It is applicable ... far apart." (plain run)
$t1a = "client-server model. "
Replace-Text $t1a $t1a

$t1b = "It is applicable when the client and server are both in the same building (and belong to the same company), but also when they are far apart."
Replace-Text $t1b $t1b

# 2) "Analysing the client-server model in detail..." paragraph -> single run
$t2 = "Analysing the client-server model in detail, we see that two processes are involved, one on the client machine and one on the server machine. "
$t2 += "Communication takes the form of the client process sending a message over the network to the server process. The client process then waits for a reply message. "
$t2 += "When the server process gets the request, it performs the requested work or looks up the requested data and sends back a reply."
Replace-Text $t2 $t2

# 3) PANs/Bluetooth paragraph merge
$t3 = ") let devices communicate over the range of a person. A common example is a wireless network that connects a computer with its peripherals. Bluetooth, a short-range wireless network"
Replace-Text $t3 $t3

# 4) "1.2.2 Local Area Network"
$t4 = "1.2.2 Local Area Network"
Replace-Text $t4 $t4

# 5) "1.2.3 Metropolitan Area Network"
$t5 = "1.2.3 Metropolitan Area Network"
Replace-Text $t5 $t5

# 6) "1.2.4 Wide Area Networks"
$t6 = "1.2.4 Wide Area Networks"
Replace-Text $t6 $t6

# 7) "Basically, a protocol ... impossible."
$t7 = "Basically, a protocol is an agreement between the communicating parties on how communication is to proceed. "
$t7 += "Violating the protocol will make communication more difficult, if not completely impossible."
Replace-Text $t7 $t7

# 8) "clearcut" merge (also removes the proofErr spellStart/spellEnd wrapper)
$t8 = " a specific collection of well-understood functions. In addition to minimizing the amount of information that must be passed between layers, clearcut interfaces also make it simpler to replace one layer with a completely different protocol or implementation"
Replace-Text $t8 $t8

# 9) Circuit switching caption merge
$t9 = " - In this example, each office has three incoming lines and three outgoing lines. When a call passes through a switching office, a physical connection is (conceptually) established between the line on which the call came in and one of the output lines, as shown by the dotted lines. "
Replace-Text $t9 $t9

# 10) Packet switching caption merge
$t10 = " - With this technology, packets are sent as soon as they are available. There is no need to set up a dedicated path in advance, unlike with circuit switching. It is up to routers to use store-and-forward transmission to send each packet on its way to the destination on its own. This procedure is unlike circuit switching, in which the result of the connection setup is the reservation of bandwidth all the way from the sender to the receiver. "
Replace-Text $t10 $t10

# 11) "With packet switching ... packets are sent)." merge
$t11 = "With packet switching there is no fixed path, so different packets can follow different paths, depending on network conditions at the time they are sent, and they may arrive out of order. "
$t11 += "Because no bandwidth is reserved with packet switching, packets may have to wait to be forwarded. This introduces queuing delay and congestion if many packets are sent at the same time. "
$t11 += "On the other hand, there is no danger of getting a busy signal and being unable to use the network. Thus, congestion occurs at different times with circuit switching (at setup time) and packet switching (when packets are sent)."
Replace-Text $t11 $t11

Write-Output "done with text merges"
